# The author deleted two data rows from Sheet1 (the rows for genes "LTV1"
# and "SRO9"), which removes their now-unused shared-string entries too and
# shifts every following row up. Locate each row by its label (column A)
# rather than a hard-coded row number so the script is robust to exactly
# how the rows are currently laid out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row labeled "LTV1" (originally row 252).
$ltv1 = $ws.Cells.Find("LTV1")
$ltv1.EntireRow.Delete()

# Delete the entire row labeled "SRO9" (originally row 259, now row 258
# after the previous deletion shifted things up by one).
$sro9 = $ws.Cells.Find("SRO9")
$sro9.EntireRow.Delete()

# Excel leaves the selection sitting just past the old used range after
# the rows above it were removed.
$ws.Range("F281").Select()

# Portrait page orientation was (re)applied on save.
$ws.PageSetup.Orientation = 1
